$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3247571
$ws.Range("I33").Value = 876.8077
$ws.Range("K33").Value = 876.8077
$ws.Range("M33").Value = -647.8077

$ws.Range("H116").Value = 7944.9375
$ws.Range("I116").Value = 9959.916999999999
$ws.Range("J116").Value = 1900
$ws.Range("K116").Value = 9959.916999999999
$ws.Range("L116").Value = 1900
$ws.Range("M116").Value = -6517.916999999999
$ws.Range("N116").Value = -8784

$ws.Range("H127").Value = 2071.2222
$ws.Range("I127").Value = 694.25
$ws.Range("J127").Value = 2243.3438
$ws.Range("K127").Value = 2082.75
$ws.Range("L127").Value = 6730.0314
$ws.Range("M127").Value = 2877.25
$ws.Range("N127").Value = -16650.0314

$ws.Range("H132").Value = 896.34485
$ws.Range("I132").Value = 833.1111
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 2499.3333
$ws.Range("L132").Value = 5250
$ws.Range("M132").Value = 30.66670000000022
$ws.Range("N132").Value = -10310

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2920.6
$ws.Range("I88").Value = 2798
$ws.Range("K88").Value = 2798
$ws.Range("M88").Value = -2392

$ws.Range("H91").Value = 2920.6
$ws.Range("I91").Value = 2798
$ws.Range("K91").Value = 2798
$ws.Range("M91").Value = -1394

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7938294.5
$ws.Range("I86").Value = 12347373
$ws.Range("J86").Value = 1954.0667
$ws.Range("K86").Value = 12347373
$ws.Range("L86").Value = 1954.0667
$ws.Range("M86").Value = -12346250
$ws.Range("N86").Value = -4200.0667

$ws.Range("H89").Value = 7938294.5
$ws.Range("I89").Value = 12347373
$ws.Range("J89").Value = 1954.0667
$ws.Range("K89").Value = 61736865
$ws.Range("L89").Value = 9770.333500000001
$ws.Range("M89").Value = -61731249
$ws.Range("N89").Value = -21002.3335

$ws.Range("H99").Value = 142858350
$ws.Range("I99").Value = 500000100
$ws.Range("J99").Value = 1649.6
$ws.Range("K99").Value = 500000100
$ws.Range("L99").Value = 1649.6
$ws.Range("M99").Value = -499998602
$ws.Range("N99").Value = -4645.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3001.976
$ws.Range("I31").Value = 1040.85
$ws.Range("K31").Value = 1040.85
$ws.Range("M31").Value = -745.8499999999999

$ws.Range("H34").Value = 3001.976
$ws.Range("I34").Value = 1040.85
$ws.Range("K34").Value = 1040.85
$ws.Range("M34").Value = -838.8499999999999

$ws.Range("H36").Value = 14100
$ws.Range("J36").Value = 14100
$ws.Range("L36").Value = 14100
$ws.Range("N36").Value = -14876

$ws.Range("H40").Value = 14100
$ws.Range("J40").Value = 14100
$ws.Range("L40").Value = 14100
$ws.Range("N40").Value = -14420

$ws.Range("H58").Value = 1139.1063
$ws.Range("I58").Value = 791.375
$ws.Range("J58").Value = 1880.9333
$ws.Range("K58").Value = 791.375
$ws.Range("L58").Value = 1880.9333
$ws.Range("M58").Value = -588.375
$ws.Range("N58").Value = -2286.9333

$ws.Range("H62").Value = 12500
$ws.Range("I62").Value = 20000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 20000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -19376
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 12500
$ws.Range("I65").Value = 20000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 100000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -96880
$ws.Range("N65").Value = -31240

$ws.Range("H122").Value = 1016.38464
$ws.Range("I122").Value = 762.25
$ws.Range("J122").Value = 1423
$ws.Range("K122").Value = 2286.75
$ws.Range("L122").Value = 4269
$ws.Range("M122").Value = 163.25
$ws.Range("N122").Value = -9169

$ws.Range("H132").Value = 1906
$ws.Range("I132").Value = 1656.25
$ws.Range("J132").Value = 2632.5454
$ws.Range("K132").Value = 4968.75
$ws.Range("L132").Value = 7897.6362
$ws.Range("M132").Value = -2438.75
$ws.Range("N132").Value = -12957.6362

$ws.Range("H134").Value = 2913.5386
$ws.Range("I134").Value = 4054.9333
$ws.Range("J134").Value = 1357.091
$ws.Range("K134").Value = 12164.7999
$ws.Range("L134").Value = 4071.273
$ws.Range("M134").Value = -9629.7999
$ws.Range("N134").Value = -9141.272999999999

$ws.Range("H136").Value = 1139.1063
$ws.Range("I136").Value = 791.375
$ws.Range("J136").Value = 1880.9333
$ws.Range("K136").Value = 2374.125
$ws.Range("L136").Value = 5642.7999
$ws.Range("M136").Value = 175.875
$ws.Range("N136").Value = -10742.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 6000529
$ws.Range("I113").Value = 12500447
$ws.Range("J113").Value = 1667250
$ws.Range("K113").Value = 37501341
$ws.Range("L113").Value = 5001750
$ws.Range("M113").Value = -37499171
$ws.Range("N113").Value = -5006090

$ws.Range("H131").Value = 1205685.1
$ws.Range("I131").Value = 5263470.5
$ws.Range("J131").Value = 1030.1875
$ws.Range("K131").Value = 15790411.5
$ws.Range("L131").Value = 3090.5625
$ws.Range("M131").Value = -15785371.5
$ws.Range("N131").Value = -13170.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H43").Value = 5001000
$ws.Range("I43").Value = 10000000
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 10000000
$ws.Range("L43").Value = 2000
$ws.Range("M43").Value = -9999849
$ws.Range("N43").Value = -2302

$ws.Range("H46").Value = 11210.25
$ws.Range("I46").Value = 9920.5
$ws.Range("J46").Value = 12500
$ws.Range("K46").Value = 9920.5
$ws.Range("L46").Value = 12500
$ws.Range("M46").Value = -9764.5
$ws.Range("N46").Value = -12812

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H57").Value = 9422
$ws.Range("J57").Value = 9422
$ws.Range("L57").Value = 9422
$ws.Range("N57").Value = -11062

$ws.Range("H70").Value = 5655.4814
$ws.Range("I70").Value = 5608.591
$ws.Range("J70").Value = 5861.8
$ws.Range("K70").Value = 5608.591
$ws.Range("L70").Value = 5861.8
$ws.Range("M70").Value = -5338.591
$ws.Range("N70").Value = -6401.8

$ws.Range("H73").Value = 5655.4814
$ws.Range("I73").Value = 5608.591
$ws.Range("J73").Value = 5861.8
$ws.Range("K73").Value = 5608.591
$ws.Range("L73").Value = 5861.8
$ws.Range("M73").Value = -4672.591
$ws.Range("N73").Value = -7733.8
